# New crime data collected - weekly CompStat update (19th Precinct)
# Update the report header (week/volume number + date range), and refresh
# the crime-complaint figures for the week-to-date / 28-day / year-to-date
# table (rows 14-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 31   Number  36" -> "...  37"
#              "Report Covering the Week  9/2/2024  Through  9/8/2024"
#              -> "...  9/9/2024  Through  9/15/2024"
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/9/2024  Through  9/15/2024"

# ---------------------------------------------------------------------
# Helper: some cells were previously text placeholders ("0" / "***.*")
# and must become real numbers. Copying the NumberFormat from a sibling
# numeric cell in the same row flips the cell's style from the text
# style (s=14) to the numeric style (s=15/16) to match.
# ---------------------------------------------------------------------
function Set-NumericCell($cell, $value, $formatSourceCell) {
    $cell.NumberFormat = $formatSourceCell.NumberFormat
    $cell.Value = $value
}

# Row 14 - Murder
$ws.Range("N14").Value = -80

# Row 15 - Rape
$ws.Range("M15").Value = -10

# Row 16 - Robbery
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -28.571428571428
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 31.25
$ws.Range("I16").Value = 149
$ws.Range("J16").Value = 146
$ws.Range("K16").Value = 2.054794520547
$ws.Range("L16").Value = -14.857142857142
$ws.Range("M16").Value = 53.608247422680
$ws.Range("N16").Value = -84.479166666666

# Row 17 - Fel. Assault (C17 was a text placeholder "0", now a real number)
Set-NumericCell $ws.Range("C17") 2 $ws.Range("D17")
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -37.5
$ws.Range("I17").Value = 125
$ws.Range("J17").Value = 138
$ws.Range("K17").Value = -9.420289855072
$ws.Range("L17").Value = -3.100775193798
$ws.Range("M17").Value = 89.393939393939
$ws.Range("N17").Value = -40.191387559808

# Row 18 - Burglary
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -71.428571428571
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -65.217391304347
$ws.Range("I18").Value = 160
$ws.Range("J18").Value = 186
$ws.Range("K18").Value = -13.978494623655
$ws.Range("L18").Value = -10.112359550561
$ws.Range("M18").Value = -11.111111111111
$ws.Range("N18").Value = -91.894630192502

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 40
$ws.Range("D19").Value = 38
$ws.Range("E19").Value = 5.263157894736
$ws.Range("F19").Value = 142
$ws.Range("G19").Value = 146
$ws.Range("H19").Value = -2.739726027397
$ws.Range("I19").Value = 1125
$ws.Range("J19").Value = 1219
$ws.Range("K19").Value = -7.711238720262
$ws.Range("L19").Value = -8.013082583810
$ws.Range("M19").Value = 27.695800227014
$ws.Range("N19").Value = -56.293706293706

# Row 20 - G.L.A.
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -83.333333333333
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -37.5
$ws.Range("I20").Value = 67
$ws.Range("J20").Value = 123
$ws.Range("K20").Value = -45.528455284552
$ws.Range("L20").Value = -50
$ws.Range("M20").Value = 21.818181818181
$ws.Range("N20").Value = -97.201336675020

# Row 21 - TOTAL
$ws.Range("C21").Value = 50
$ws.Range("D21").Value = 64
$ws.Range("E21").Value = -21.875
$ws.Range("F21").Value = 192
$ws.Range("G21").Value = 217
$ws.Range("H21").Value = -11.520737327188
$ws.Range("I21").Value = 1636
$ws.Range("J21").Value = 1823
$ws.Range("K21").Value = -10.257816785518
$ws.Range("L21").Value = -11.567567567567
$ws.Range("M21").Value = 26.821705426356
$ws.Range("N21").Value = -79.891838741396

# Row 22 - Transit (D22/E22 were text placeholders, now real numbers)
$ws.Range("C22").Value = 1
Set-NumericCell $ws.Range("D22") 2 $ws.Range("C22")
Set-NumericCell $ws.Range("E22") -50 $ws.Range("E18")
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 32
$ws.Range("J22").Value = 34
$ws.Range("K22").Value = -5.882352941176
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 52.380952380952
# N22 remains the "***.*" text placeholder (unchanged)

# Row 23 - Housing (D23/E23 were text placeholders, now real numbers)
$ws.Range("C23").Value = 1
Set-NumericCell $ws.Range("D23") 2 $ws.Range("C23")
Set-NumericCell $ws.Range("E23") -50 $ws.Range("E18")
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 21
$ws.Range("J23").Value = 18
$ws.Range("K23").Value = 16.666666666666
$ws.Range("L23").Value = -12.5
$ws.Range("M23").Value = 10.526315789473
# N23 remains the "***.*" text placeholder (unchanged)

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 88
$ws.Range("D24").Value = 59
$ws.Range("E24").Value = 49.152542372881
$ws.Range("F24").Value = 291
$ws.Range("G24").Value = 234
$ws.Range("H24").Value = 24.358974358974
$ws.Range("I24").Value = 2280
$ws.Range("J24").Value = 2302
$ws.Range("K24").Value = -0.955690703735
$ws.Range("L24").Value = -18.338108882521
$ws.Range("M24").Value = 92.893401015228
# N24 remains the "***.*" text placeholder (unchanged)

# Row 25 - Retail Theft
$ws.Range("C25").Value = 72
$ws.Range("D25").Value = 47
$ws.Range("E25").Value = 53.191489361702
$ws.Range("F25").Value = 247
$ws.Range("G25").Value = 189
$ws.Range("H25").Value = 30.687830687830
$ws.Range("I25").Value = 1991
$ws.Range("J25").Value = 2031
$ws.Range("K25").Value = -1.969473165928
$ws.Range("L25").Value = -23.628691983122
# M25/N25 remain the "***.*" text placeholders (unchanged)

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 14.285714285714
$ws.Range("F26").Value = 28
$ws.Range("G26").Value = 15
$ws.Range("H26").Value = 86.666666666666
$ws.Range("I26").Value = 261
$ws.Range("J26").Value = 224
$ws.Range("K26").Value = 16.517857142857
$ws.Range("L26").Value = -2.973977695167
$ws.Range("M26").Value = 5.668016194331
# N26 remains the "***.*" text placeholder (unchanged)

# Row 28 - Other Sex Crimes
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 300
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 83
$ws.Range("J28").Value = 67
$ws.Range("K28").Value = 23.880597014925
$ws.Range("L28").Value = 33.870967741935
# M28/N28 remain the "***.*" text placeholders (unchanged)

# Row 31 - Hate Crimes
$ws.Range("G31").Value = 1
